$wb = $excel.ActiveWorkbook

# "focal" sheet (sheetId 2) gets new columns I and J with header/data,
# and becomes the active sheet/tab.
$ws = $wb.Worksheets.Item("focal")

$ws.Range("I1").Value = "y"
$ws.Range("J1").Value = "x"

$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 23.2

$ws.Range("I3").Value = 7.7
$ws.Range("J3").Value = 18

$ws.Range("I4").Value = 6.5
$ws.Range("J4").Value = 13.8

$ws.Range("I5").Value = 10.6
$ws.Range("J5").Value = 8.8

# Make this sheet the active tab and select I6, matching the new view state.
$ws.Activate()
$ws.Range("I6").Select()
